$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Sending cluster","Ligand symbol","Receptor symbol","Target cluster","Ligand-expressing cells","Ligand detection rate","Ligand average expression value","Ligand total expression value","Ligand derived specificity of average expression value","Ligand derived specificity of total expression value","Receptor-expressing cells","Receptor detection rate","Receptor average expression value","Receptor total expression value","Receptor derived specificity of average expression value","Receptor derived specificity of total expression value","Edge average expression weight","Edge total expression weight","Edge average expression derived specificity","Edge total expression derived specificity")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col5a3"
$ws.Range("C2").Value = "Sdc3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.894567
$ws.Range("H2").Value = 2.683701
$ws.Range("I2").Value = 0.007903265526675987
$ws.Range("J2").Value = 0.007903265526675987
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.12444933333333
$ws.Range("N2").Value = 78.373348
$ws.Range("O2").Value = 0.7238861157526749
$ws.Range("P2").Value = 0.7238861157526749
$ws.Range("Q2").Value = 23.370070266772
$ws.Range("R2").Value = 210.330632400948
$ws.Range("S2").Value = 0.005721064183867499
$ws.Range("T2").Value = 0.005721064183867499

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col5a3"
$ws.Range("C3").Value = "Sdc3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.894567
$ws.Range("H3").Value = 2.683701
$ws.Range("I3").Value = 0.007903265526675987
$ws.Range("J3").Value = 0.007903265526675987
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.818542
$ws.Range("N3").Value = 11.455626
$ws.Range("O3").Value = 0.1058085282850919
$ws.Range("P3").Value = 0.1058085282850919
$ws.Range("Q3").Value = 3.415941661314
$ws.Range("R3").Value = 30.743474951826
$ws.Range("S3").Value = 0.0008362328940238882
$ws.Range("T3").Value = 0.0008362328940238882

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col5a3"
$ws.Range("C4").Value = "Sdc3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.894567
$ws.Range("H4").Value = 2.683701
$ws.Range("I4").Value = 0.007903265526675987
$ws.Range("J4").Value = 0.007903265526675987
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.146179
$ws.Range("N4").Value = 18.438537
$ws.Range("O4").Value = 0.1703053559622332
$ws.Range("P4").Value = 0.1703053559622332
$ws.Range("Q4").Value = 5.498168909493
$ws.Range("R4").Value = 49.483520185437
$ws.Range("S4").Value = 0.001345968448784601
$ws.Range("T4").Value = 0.001345968448784601

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col5a3"
$ws.Range("C5").Value = "Sdc3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 97.363968
$ws.Range("H5").Value = 292.091904
$ws.Range("I5").Value = 0.860185197793775
$ws.Range("J5").Value = 0.860185197793775
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 26.12444933333333
$ws.Range("N5").Value = 78.373348
$ws.Range("O5").Value = 0.7238861157526749
$ws.Range("P5").Value = 0.7238861157526749
$ws.Range("Q5").Value = 2543.580048908288
$ws.Range("R5").Value = 22892.22044017459
$ws.Range("S5").Value = 0.6226761216588822
$ws.Range("T5").Value = 0.6226761216588822

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col5a3"
$ws.Range("C6").Value = "Sdc3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 97.363968
$ws.Range("H6").Value = 292.091904
$ws.Range("I6").Value = 0.860185197793775
$ws.Range("J6").Value = 0.860185197793775
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.818542
$ws.Range("N6").Value = 11.455626
$ws.Range("O6").Value = 0.1058085282850919
$ws.Range("P6").Value = 0.1058085282850919
$ws.Range("Q6").Value = 371.788401094656
$ws.Range("R6").Value = 3346.095609851904
$ws.Range("S6").Value = 0.09101492983118005
$ws.Range("T6").Value = 0.09101492983118005

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col5a3"
$ws.Range("C7").Value = "Sdc3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 97.363968
$ws.Range("H7").Value = 292.091904
$ws.Range("I7").Value = 0.860185197793775
$ws.Range("J7").Value = 0.860185197793775
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.146179
$ws.Range("N7").Value = 18.438537
$ws.Range("O7").Value = 0.1703053559622332
$ws.Range("P7").Value = 0.1703053559622332
$ws.Range("Q7").Value = 598.416375478272
$ws.Range("R7").Value = 5385.747379304448
$ws.Range("S7").Value = 0.1464941463037129
$ws.Range("T7").Value = 0.1464941463037129

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Col5a3"
$ws.Range("C8").Value = "Sdc3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 14.93100633333333
$ws.Range("H8").Value = 44.793019
$ws.Range("I8").Value = 0.1319115366795491
$ws.Range("J8").Value = 0.1319115366795491
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 26.12444933333333
$ws.Range("N8").Value = 78.373348
$ws.Range("O8").Value = 0.7238861157526749
$ws.Range("P8").Value = 0.7238861157526749
$ws.Range("Q8").Value = 390.0643184508457
$ws.Range("R8").Value = 3510.578866057612
$ws.Range("S8").Value = 0.09548892990992529
$ws.Range("T8").Value = 0.09548892990992527

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Col5a3"
$ws.Range("C9").Value = "Sdc3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 14.93100633333333
$ws.Range("H9").Value = 44.793019
$ws.Range("I9").Value = 0.1319115366795491
$ws.Range("J9").Value = 0.1319115366795491
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.818542
$ws.Range("N9").Value = 11.455626
$ws.Range("O9").Value = 0.1058085282850919
$ws.Range("P9").Value = 0.1058085282850919
$ws.Range("Q9").Value = 57.01467478609934
$ws.Range("R9").Value = 513.132073074894
$ws.Range("S9").Value = 0.01395736555988801
$ws.Range("T9").Value = 0.01395736555988801

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Col5a3"
$ws.Range("C10").Value = "Sdc3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 14.93100633333333
$ws.Range("H10").Value = 44.793019
$ws.Range("I10").Value = 0.1319115366795491
$ws.Range("J10").Value = 0.1319115366795491
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 6.146179
$ws.Range("N10").Value = 18.438537
$ws.Range("O10").Value = 0.1703053559622332
$ws.Range("P10").Value = 0.1703053559622332
$ws.Range("Q10").Value = 91.76863757480034
$ws.Range("R10").Value = 825.917738173203
$ws.Range("S10").Value = 0.02246524120973579
$ws.Range("T10").Value = 0.02246524120973579

